$wb = $excel.ActiveWorkbook

$wsOrdner = $wb.Worksheets.Item("Ordner")
$wsSeiten = $wb.Worksheets.Item("Seiten")

# --- Update "In Ordner" (column B) values: increment by 1 for data rows 2-5 ---
$wsSeiten.Range("B2").Value = 2
$wsSeiten.Range("B3").Value = 2
$wsSeiten.Range("B4").Value = 3
$wsSeiten.Range("B5").Value = 3

# --- Update D4 text: append a new line with alt+enter note ---
$wsSeiten.Range("D4").Value = "Es ist auch möglich, normalen Text zu verwenden. Wir werden das Beste daraus machen.`nVerwenden Sie die Tastenkombination alt+enter in Excel, um Zeilenumbrüche hinzuzufügen."
$wsSeiten.Range("D4").WrapText = $true
$wsSeiten.Rows.Item(4).RowHeight = 28.8

# --- Autofit columns on both sheets ---
$wsOrdner.Columns.AutoFit() | Out-Null
$wsSeiten.Columns.AutoFit() | Out-Null

# --- Update selections (cosmetic cursor position) ---
$wsOrdner.Range("E7").Select()
$wsSeiten.Range("G9").Select()
